$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Remove the obsolete "user_id" column (column A); the remaining columns
# (user_name, user_login, Senha, user_pass, user_level) shift left to A:E.
$ws.Columns.Item(1).Delete() | Out-Null

# Restore the cursor/selection to where the author left it after editing.
$ws.Range("I6").Select() | Out-Null
